$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Plain text cells (coin names / URLs) - safe to assign directly
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

# Numeric-looking text cells (price/volume/hour) - force Text format to preserve exact string (leading/trailing zeros, % sign, no float coercion)
$c = $ws.Range("D2"); $c.NumberFormat = "@"; $c.Value = '258.50'
$c = $ws.Range("E2"); $c.NumberFormat = "@"; $c.Value = '1.35%'
$c = $ws.Range("G2"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D3"); $c.NumberFormat = "@"; $c.Value = '26.93'
$c = $ws.Range("E3"); $c.NumberFormat = "@"; $c.Value = '-3.90%'
$c = $ws.Range("G3"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D4"); $c.NumberFormat = "@"; $c.Value = '4.642'
$c = $ws.Range("E4"); $c.NumberFormat = "@"; $c.Value = '-13.20%'
$c = $ws.Range("G4"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D5"); $c.NumberFormat = "@"; $c.Value = '0.05973'
$c = $ws.Range("G5"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("E6"); $c.NumberFormat = "@"; $c.Value = '-0.36%'
$c = $ws.Range("G6"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D7"); $c.NumberFormat = "@"; $c.Value = '0.8754'
$c = $ws.Range("E7"); $c.NumberFormat = "@"; $c.Value = '1.44%'
$c = $ws.Range("G7"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D8"); $c.NumberFormat = "@"; $c.Value = '0.9553'
$c = $ws.Range("E8"); $c.NumberFormat = "@"; $c.Value = '4.10%'
$c = $ws.Range("G8"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D9"); $c.NumberFormat = "@"; $c.Value = '0.0006089'
$c = $ws.Range("E9"); $c.NumberFormat = "@"; $c.Value = '-94.27%'
$c = $ws.Range("G9"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D10"); $c.NumberFormat = "@"; $c.Value = '0.1416'
$c = $ws.Range("E10"); $c.NumberFormat = "@"; $c.Value = '-0.25%'
$c = $ws.Range("G10"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D11"); $c.NumberFormat = "@"; $c.Value = '0.07188'
$c = $ws.Range("E11"); $c.NumberFormat = "@"; $c.Value = '0.18%'
$c = $ws.Range("G11"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D12"); $c.NumberFormat = "@"; $c.Value = '0.03137'
$c = $ws.Range("E12"); $c.NumberFormat = "@"; $c.Value = '-2.10%'
$c = $ws.Range("G12"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D13"); $c.NumberFormat = "@"; $c.Value = '0.09234'
$c = $ws.Range("E13"); $c.NumberFormat = "@"; $c.Value = '-0.01%'
$c = $ws.Range("G13"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D14"); $c.NumberFormat = "@"; $c.Value = '0.001543'
$c = $ws.Range("E14"); $c.NumberFormat = "@"; $c.Value = '-0.38%'
$c = $ws.Range("G14"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D15"); $c.NumberFormat = "@"; $c.Value = '0.005953'
$c = $ws.Range("E15"); $c.NumberFormat = "@"; $c.Value = '1.32%'
$c = $ws.Range("G15"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D16"); $c.NumberFormat = "@"; $c.Value = '3.487'
$c = $ws.Range("E16"); $c.NumberFormat = "@"; $c.Value = '-0.42%'
$c = $ws.Range("G16"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D17"); $c.NumberFormat = "@"; $c.Value = '3.210'
$c = $ws.Range("E17"); $c.NumberFormat = "@"; $c.Value = '-0.43%'
$c = $ws.Range("G17"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D18"); $c.NumberFormat = "@"; $c.Value = '2.219'
$c = $ws.Range("E18"); $c.NumberFormat = "@"; $c.Value = '-1.47%'
$c = $ws.Range("G18"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D19"); $c.NumberFormat = "@"; $c.Value = '0.3136'
$c = $ws.Range("E19"); $c.NumberFormat = "@"; $c.Value = '-1.00%'
$c = $ws.Range("G19"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D20"); $c.NumberFormat = "@"; $c.Value = '0.03608'
$c = $ws.Range("E20"); $c.NumberFormat = "@"; $c.Value = '4.71%'
$c = $ws.Range("G20"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("E21"); $c.NumberFormat = "@"; $c.Value = '-0.97%'
$c = $ws.Range("G21"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D22"); $c.NumberFormat = "@"; $c.Value = '3.533'
$c = $ws.Range("E22"); $c.NumberFormat = "@"; $c.Value = '0.14%'
$c = $ws.Range("G22"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D23"); $c.NumberFormat = "@"; $c.Value = '0.04221'
$c = $ws.Range("E23"); $c.NumberFormat = "@"; $c.Value = '1.81%'
$c = $ws.Range("G23"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("E24"); $c.NumberFormat = "@"; $c.Value = '0.05%'
$c = $ws.Range("G24"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D25"); $c.NumberFormat = "@"; $c.Value = '0.001220'
$c = $ws.Range("E25"); $c.NumberFormat = "@"; $c.Value = '-0.20%'
$c = $ws.Range("G25"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("E26"); $c.NumberFormat = "@"; $c.Value = '-11.72%'
$c = $ws.Range("G26"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("E27"); $c.NumberFormat = "@"; $c.Value = '0.03%'
$c = $ws.Range("G27"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("E28"); $c.NumberFormat = "@"; $c.Value = '-22.96%'
$c = $ws.Range("G28"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("G29"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("G30"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("G31"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("G32"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("G33"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("G34"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("G35"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("G36"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("G37"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("G38"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("G39"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D40"); $c.NumberFormat = "@"; $c.Value = '0.03849'
$c = $ws.Range("E40"); $c.NumberFormat = "@"; $c.Value = '0.09%'
$c = $ws.Range("G40"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D41"); $c.NumberFormat = "@"; $c.Value = '0.006029'
$c = $ws.Range("E41"); $c.NumberFormat = "@"; $c.Value = '58.88%'
$c = $ws.Range("G41"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D42"); $c.NumberFormat = "@"; $c.Value = '0.1104'
$c = $ws.Range("E42"); $c.NumberFormat = "@"; $c.Value = '0.41%'
$c = $ws.Range("G42"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D43"); $c.NumberFormat = "@"; $c.Value = '0.002200'
$c = $ws.Range("E43"); $c.NumberFormat = "@"; $c.Value = '-8.31%'
$c = $ws.Range("G43"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("E44"); $c.NumberFormat = "@"; $c.Value = '10.69%'
$c = $ws.Range("G44"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D45"); $c.NumberFormat = "@"; $c.Value = '0.00005495'
$c = $ws.Range("E45"); $c.NumberFormat = "@"; $c.Value = '3.95%'
$c = $ws.Range("G45"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("E46"); $c.NumberFormat = "@"; $c.Value = '0.00%'
$c = $ws.Range("G46"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("E47"); $c.NumberFormat = "@"; $c.Value = '-14.49%'
$c = $ws.Range("G47"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D48"); $c.NumberFormat = "@"; $c.Value = '0.002126'
$c = $ws.Range("E48"); $c.NumberFormat = "@"; $c.Value = '-3.79%'
$c = $ws.Range("G48"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D49"); $c.NumberFormat = "@"; $c.Value = '0.00002100'
$c = $ws.Range("E49"); $c.NumberFormat = "@"; $c.Value = '0.00%'
$c = $ws.Range("G49"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("D50"); $c.NumberFormat = "@"; $c.Value = '0.0002000'
$c = $ws.Range("E50"); $c.NumberFormat = "@"; $c.Value = '0.00%'
$c = $ws.Range("G50"); $c.NumberFormat = "@"; $c.Value = '6'
$c = $ws.Range("G51"); $c.NumberFormat = "@"; $c.Value = '6'

Write-Output "Applied 132 cell changes"